# Update the "Förändrad" (Changed) date column for rows 2-5 from
# 2023-09-05 (serial 45174) to 2023-09-06 (serial 45175).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2").Value = 45175
$ws.Range("C3").Value = 45175
$ws.Range("C4").Value = 45175
$ws.Range("C5").Value = 45175
